$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 105
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 415
$ws.Range("I13").Value = 415
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 415
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -246
$ws.Range("N13").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1623.3334
$ws.Range("I38").Value = 348
$ws.Range("J38").Value = 8000
$ws.Range("K38").Value = 1044
$ws.Range("L38").Value = 24000
$ws.Range("M38").Value = -672
$ws.Range("N38").Value = -24744

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 642.375
$ws.Range("I55").Value = 609.1111
$ws.Range("J55").Value = 685.1429000000001
$ws.Range("K55").Value = 609.1111
$ws.Range("L55").Value = 685.1429000000001
$ws.Range("M55").Value = -395.1111
$ws.Range("N55").Value = -1113.1429

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 1113
$ws.Range("I69").Value = 1113
$ws.Range("K69").Value = 3339
$ws.Range("M69").Value = -2465

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 1113
$ws.Range("I72").Value = 1113
$ws.Range("K72").Value = 10017
$ws.Range("M72").Value = -5649

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1878.4706
$ws.Range("I112").Value = 1037.25
$ws.Range("J112").Value = 2137.3076
$ws.Range("K112").Value = 3111.75
$ws.Range("L112").Value = 6411.9228
$ws.Range("M112").Value = -2003.75
$ws.Range("N112").Value = -8627.9228

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1262.8572
$ws.Range("I135").Value = 834.5454999999999
$ws.Range("K135").Value = 7510.9095
$ws.Range("M135").Value = -4975.9095

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1198.2084
$ws.Range("I137").Value = 1042.85
$ws.Range("K137").Value = 3128.55
$ws.Range("M137").Value = -578.5499999999997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2291.8572
$ws.Range("I61").Value = 2291.8572
$ws.Range("K61").Value = 2291.8572
$ws.Range("M61").Value = -2079.8572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1090.8572
$ws.Range("I97").Value = 1087.4
$ws.Range("J97").Value = 1099.5
$ws.Range("K97").Value = 1087.4
$ws.Range("L97").Value = 1099.5
$ws.Range("M97").Value = -591.4000000000001
$ws.Range("N97").Value = -2091.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4108.7144
$ws.Range("I122").Value = 2920.3333
$ws.Range("K122").Value = 8760.999899999999
$ws.Range("M122").Value = -6310.999899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1075.8334
$ws.Range("I132").Value = 1075.8334
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3227.5002
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -697.5001999999999
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2291.8572
$ws.Range("I136").Value = 2291.8572
$ws.Range("K136").Value = 6875.571599999999
$ws.Range("M136").Value = -4325.571599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 437.33334
$ws.Range("I22").Value = 324.8
$ws.Range("K22").Value = 324.8
$ws.Range("M22").Value = -151.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 5511.7144
$ws.Range("J61").Value = 8000
$ws.Range("L61").Value = 8000
$ws.Range("N61").Value = -8626

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2573.111
$ws.Range("I94").Value = 2519.75
$ws.Range("K94").Value = 2519.75
$ws.Range("M94").Value = -2068.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 14209.5
$ws.Range("J100").Value = 14209.5
$ws.Range("L100").Value = 14209.5
$ws.Range("N100").Value = -16373.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1420.5454
$ws.Range("I105").Value = 1075.1666
$ws.Range("J105").Value = 2974.75
$ws.Range("K105").Value = 1075.1666
$ws.Range("L105").Value = 2974.75
$ws.Range("M105").Value = 671.8334
$ws.Range("N105").Value = -6468.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1633.5952
$ws.Range("I134").Value = 893.3214
$ws.Range("K134").Value = 2679.9642
$ws.Range("M134").Value = -144.9642000000003

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 127.46667
$ws.Range("I7").Value = 47
$ws.Range("K7").Value = 47
$ws.Range("M7").Value = 66

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H53").Value = 66999
$ws.Range("J53").Value = 66999
$ws.Range("L53").Value = 66999
$ws.Range("N53").Value = -68213

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2206.7144
$ws.Range("I58").Value = 1134.2858
$ws.Range("K58").Value = 1134.2858
$ws.Range("M58").Value = -931.2858000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 728
$ws.Range("I94").Value = 704
$ws.Range("K94").Value = 704
$ws.Range("M94").Value = -253

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H97").Value = 70000
$ws.Range("I97").Value = 80000
$ws.Range("J97").Value = 60000
$ws.Range("K97").Value = 80000
$ws.Range("L97").Value = 60000
$ws.Range("M97").Value = -79009
$ws.Range("N97").Value = -61982

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 4706.393
$ws.Range("J105").Value = 6905.4443
$ws.Range("L105").Value = 6905.4443
$ws.Range("N105").Value = -10399.4443

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H106").Value = 12345
$ws.Range("J106").Value = 12345
$ws.Range("L106").Value = 12345
$ws.Range("N106").Value = -14869

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3330
$ws.Range("I134").Value = 3344.4
$ws.Range("J134").Value = 3312
$ws.Range("K134").Value = 10033.2
$ws.Range("L134").Value = 9936
$ws.Range("M134").Value = -7498.200000000001
$ws.Range("N134").Value = -15006

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2206.7144
$ws.Range("I136").Value = 1134.2858
$ws.Range("K136").Value = 3402.8574
$ws.Range("M136").Value = -852.8574000000003

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 2096.2
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2096.2
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 6288.599999999999
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -6626.599999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 886.5833
$ws.Range("I33").Value = 1448.1428
$ws.Range("J33").Value = 100.4
$ws.Range("K33").Value = 8688.856800000001
$ws.Range("L33").Value = 602.4000000000001
$ws.Range("M33").Value = -8405.856800000001
$ws.Range("N33").Value = -1168.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 2500812.5
$ws.Range("J46").Value = 3334166.8
$ws.Range("L46").Value = 10002500.4
$ws.Range("N46").Value = -10002682.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 2749.75
$ws.Range("I126").Value = 2749.75
$ws.Range("K126").Value = 8249.25
$ws.Range("M126").Value = -3309.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1080.0731
$ws.Range("J131").Value = 1088.825
$ws.Range("L131").Value = 3266.475
$ws.Range("N131").Value = -13346.475

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3321
$ws.Range("I132").Value = 3812.6155
$ws.Range("J132").Value = 2522.125
$ws.Range("K132").Value = 34313.5395
$ws.Range("L132").Value = 22699.125
$ws.Range("M132").Value = -31783.5395
$ws.Range("N132").Value = -27759.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1547.6
$ws.Range("I140").Value = 1164
$ws.Range("K140").Value = 3492
$ws.Range("M140").Value = 1688

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2656.0435
$ws.Range("I102").Value = 1955.2354
$ws.Range("K102").Value = 1955.2354
$ws.Range("M102").Value = -333.2354

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1200
$ws.Range("J107").Value = 1200
$ws.Range("L107").Value = 1200
$ws.Range("N107").Value = -5040

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1189.3334
$ws.Range("I132").Value = 761.7
$ws.Range("J132").Value = 2044.6
$ws.Range("K132").Value = 2285.1
$ws.Range("L132").Value = 6133.799999999999
$ws.Range("M132").Value = 244.8999999999996
$ws.Range("N132").Value = -11193.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1389.8462
$ws.Range("I16").Value = 1255.6666
$ws.Range("K16").Value = 1255.6666
$ws.Range("M16").Value = -1085.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2333.3333
$ws.Range("I46").Value = 750
$ws.Range("K46").Value = 750
$ws.Range("M46").Value = -562

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 923.6667
$ws.Range("I55").Value = 588.2857
$ws.Range("K55").Value = 588.2857
$ws.Range("M55").Value = -415.2857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1372.7858
$ws.Range("I93").Value = 1147.1428
$ws.Range("K93").Value = 1147.1428
$ws.Range("M93").Value = 100.8571999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2259.3225
$ws.Range("I132").Value = 1882.64
$ws.Range("J132").Value = 3828.8333
$ws.Range("K132").Value = 5647.92
$ws.Range("L132").Value = 11486.4999
$ws.Range("M132").Value = -3117.92
$ws.Range("N132").Value = -16546.4999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6607.4443
$ws.Range("I136").Value = 6097.2856
$ws.Range("K136").Value = 18291.8568
$ws.Range("M136").Value = -15741.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 21633.334
$ws.Range("I32").Value = 19450
$ws.Range("J32").Value = 26000
$ws.Range("K32").Value = 19450
$ws.Range("L32").Value = 26000
$ws.Range("M32").Value = -19133
$ws.Range("N32").Value = -26634

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 10000
$ws.Range("I34").Value = 10000
$ws.Range("K34").Value = 10000
$ws.Range("M34").Value = -9797

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 59824.75
$ws.Range("J64").Value = 59824.75
$ws.Range("L64").Value = 59824.75
$ws.Range("N64").Value = -60320.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H67").Value = 59824.75
$ws.Range("J67").Value = 59824.75
$ws.Range("L67").Value = 59824.75
$ws.Range("N67").Value = -61540.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2790.3242
$ws.Range("I132").Value = 2466.8276
$ws.Range("K132").Value = 7400.4828
$ws.Range("M132").Value = -4870.4828

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1781.8
$ws.Range("I136").Value = 1020.63635
$ws.Range("K136").Value = 3061.90905
$ws.Range("M136").Value = -511.9090500000002
